# Fruta / hortaliza, semanal
# Insert a new weekly observation row for "Ají" (Vega Monumental Concepción)
# at row 105, pushing the existing rows 105-124 down to 106-125.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 105; this shifts rows 105-124
# down to 106-125 (and carries the existing row formatting, e.g. the date
# style in column D, down with them).
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new observation.
$ws.Cells.Item(105, 1).Value = 11
$ws.Cells.Item(105, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value = "Bíobío"
$ws.Cells.Item(105, 4).Value = 44782
$ws.Cells.Item(105, 5).Value = 8
$ws.Cells.Item(105, 6).Value = 100112021
$ws.Cells.Item(105, 7).Value = "Ají"
$ws.Cells.Item(105, 8).Value = "Inferno"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 40
$ws.Cells.Item(105, 11).Value = 17000
$ws.Cells.Item(105, 12).Value = 18000
$ws.Cells.Item(105, 13).Value = 17500
$ws.Cells.Item(105, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(105, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(105, 16).Value = 1458
$ws.Cells.Item(105, 17).Value = 12
$ws.Cells.Item(105, 18).Value = "Hortaliza"
